$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -13.01
$ws.Range("C4").Value = -13.457
$ws.Range("E6").Value = 12.608
$ws.Range("C7").Value = -13.422
$ws.Range("E7").Value = 12.644
$ws.Range("C8").Value = -12.585
$ws.Range("E8").Value = 12.837
$ws.Range("A11").Value = -21.803
$ws.Range("A12").Value = -21.826
$ws.Range("C12").Value = -13.232
$ws.Range("C14").Value = -12.081
$ws.Range("A15").Value = -21.178
$ws.Range("E19").Value = 12.344
$ws.Range("E21").Value = 13.204
$ws.Range("C22").Value = -13.318
$ws.Range("E24").Value = 12.623
$ws.Range("E25").Value = 12.497
